$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (rows 1-8) holds a shared CONCATENATE() formula whose
# 'mo_fecha_crea' field is stamped with TEXT(NOW(), "yyyy-mm-dd HH:mm:ss").
# The underlying data/formula didn't change - the workbook was simply
# recalculated (re-opened/re-saved) at a later moment, so every cached
# 'mo_fecha_crea' timestamp in H1:H8 needs to be refreshed to reflect
# that new evaluation time. Force a full recalculation of the workbook
# to pick up the new NOW() value everywhere the volatile formula is used.
$excel.CalculateFull()
